$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.515.89"
$ws.Range("E2").Value = "  +3.67%  "
$ws.Range("D3").Value = "1.587.13"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("D8").Value = "'24.24"
$ws.Range("E8").Value = "  +5.34%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("D13").Value = "1.580.85"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "'3.74"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "28.530.14"
$ws.Range("D17").Value = "'63.07"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "'230.84"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("D25").Value = "'152.08"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "'0.0470"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'3.17"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "1.388.15"
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("E36").Value = "  -10.42%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "'2.63"
$ws.Range("E38").Value = "  +11.12%  "
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "'0.541"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'0.812"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'1.87"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'0.981"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'62.92"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "1.723.75"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'86.95"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'0.0522"
$ws.Range("E51").Value = "  -1.19%  "
